$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: apply the "last row" border formatting (currently on row 23) to row 22 ---
$ws.Range("B23:J23").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Step 2: update the period labels (column E, rows 16-22), newest period first ---
$ws.Range("E16").Value = "2411"
$ws.Range("E17").Value = "2410"
$ws.Range("E18").Value = "2409"
$ws.Range("E19").Value = "2408"
$ws.Range("E20").Value = "2407"
$ws.Range("E21").Value = "2406"
$ws.Range("E22").Value = "2405"

# --- Step 3: update the "Valor Mora" amounts (column F) for the table rows ---
$ws.Range("F16").Value = 52000
$ws.Range("F17").Value = 52000
$ws.Range("F18").Value = 52000
$ws.Range("F19").Value = 52000
$ws.Range("F20").Value = 52000
$ws.Range("F21").Value = 52000
$ws.Range("F22").Value = 38133

# --- Step 4: remove the now-unused last data row (shifts everything below up by one row) ---
$ws.Rows(23).Delete()

# --- Step 5: swap the "Novedad de Ingreso" / "Novedad de Retiro" header columns ---
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"

# --- Step 6: update summary figures ---
$ws.Range("E11").Value = 350133
$ws.Range("F13").Value = 7

Write-Output "done"
